$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-16 Thursday" "2025-01-17 Friday"

Replace-Text "389÷8=48, 5" "933÷9=103, 6"
Replace-Text "730÷5=146, 0" "443÷9=49, 2"
Replace-Text "160÷4=40, 0" "700÷6=116, 4"
Replace-Text "727÷9=80, 7" "579÷9=64, 3"
Replace-Text "772÷4=193, 0" "866÷5=173, 1"

Replace-Text "118÷5=23, 3" "825÷6=137, 3"
Replace-Text "609÷2=304, 1" "309÷6=51, 3"
Replace-Text "406÷2=203, 0" "410÷8=51, 2"
Replace-Text "294÷8=36, 6" "154÷5=30, 4"
Replace-Text "414÷6=69, 0" "472÷8=59, 0"

Replace-Text "650÷9=72, 2" "763÷4=190, 3"
Replace-Text "214÷3=71, 1" "738÷2=369, 0"
Replace-Text "654÷9=72, 6" "355÷5=71, 0"
Replace-Text "507÷2=253, 1" "321÷8=40, 1"
Replace-Text "441÷8=55, 1" "970÷8=121, 2"

Replace-Text "574÷3=191, 1" "478÷3=159, 1"
Replace-Text "538÷6=89, 4" "956÷7=136, 4"
Replace-Text "213÷6=35, 3" "978÷8=122, 2"
Replace-Text "418÷8=52, 2" "788÷3=262, 2"
Replace-Text "533÷6=88, 5" "830÷4=207, 2"

Replace-Text "390÷3=130, 0" "782÷8=97, 6"
Replace-Text "625÷2=312, 1" "643÷9=71, 4"
Replace-Text "705÷4=176, 1" "136÷9=15, 1"
Replace-Text "859÷6=143, 1" "451÷7=64, 3"
Replace-Text "186÷8=23, 2" "878÷9=97, 5"

Write-Output "Done"
